# Update NumPapers (column D) and Categories (column E) values for rows 2-46
# on Sheet1, per the revised keyword-by-year counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @(5, 2)
    3  = @(5, 0)
    4  = @(5, 1)
    5  = @(5, 0)
    6  = @(5, 4)
    7  = @(6, 3)
    8  = @(7, 2)
    9  = @(3, 0)
    10 = @(8, 5)
    11 = @(8, 5)
    12 = @(6, 4)
    13 = @(6, 6)
    14 = @(7, 6)
    15 = @(12, 3)
    16 = @(10, 3)
    17 = @(11, 0)
    18 = @(11, 5)
    19 = @(7, 0)
    20 = @(15, 12)
    21 = @(10, 8)
    22 = @(10, 9)
    23 = @(9, 8)
    24 = @(12, 3)
    25 = @(13, 11)
    26 = @(12, 5)
    27 = @(12, 7)
    28 = @(11, 3)
    29 = @(16, 3)
    30 = @(18, 1)
    31 = @(14, 9)
    32 = @(17, 3)
    33 = @(9, 6)
    34 = @(21, 21)
    35 = @(22, 15)
    36 = @(14, 6)
    37 = @(19, 11)
    38 = @(9, 1)
    39 = @(22, 22)
    40 = @(21, 16)
    41 = @(15, 9)
    42 = @(21, 15)
    43 = @(12, 6)
    44 = @(23, 6)
    45 = @(23, 16)
    46 = @(17, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}
